$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (row 9) for the "memToReg" control signal.
$ws.Range("A9").Value = "memToReg"

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "X"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0

$ws.Range("A9:O9").HorizontalAlignment = -4108

$ws.Range("F11").Select()
